# Update the "Revision" sheet part numbers / qty to reflect the
# temperature-sensor part-list refresh described in the commit message.
#
# Row 6 (zener):        BZT52C15-FDITR-ND  -> BZT52C15S-FDICT-ND  (flag w/ yellow highlight)
# Row 7 (fuse holder):  F1498-ND           -> 486-2019-ND          (qty 1 -> 2)
# Row 8 (barrel jack):  CP-202B-ND         -> "CP-048H-ND "        (flag w/ yellow highlight)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Revision")

# Row 8: barrel jack part number changed (note trailing space), highlighted yellow.
$ws.Range("B8").Value = "CP-048H-ND "
$ws.Range("B8").Interior.Color = 65535

# Row 7: fuse holder part number + quantity changed, no highlight.
$ws.Range("B7").Value = "486-2019-ND"
$ws.Range("C7").Value = 2

# Row 6: zener part number changed, highlighted yellow like the other
# previously-revised rows (B2, B10, B11, B14).
$ws.Range("B6").Value = "BZT52C15S-FDICT-ND"
$ws.Range("B6").Interior.Color = 65535

# Move the saved selection on both sheets to B6, matching the refreshed file.
$ws1 = $wb.Worksheets.Item("First Order")
$ws1.Activate()
$ws1.Range("B6").Select()

$ws.Activate()
$ws.Range("B6").Select()
